$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new filename (column A) and new curvature value (column F, only where changed)
$updates = @(
    @{ Row = 2;  A = "01.jpg"; F = 0.02431618966192055 },
    @{ Row = 3;  A = "02.jpg"; F = 0.02228918141274923 },
    @{ Row = 4;  A = "03.jpg"; F = 0.02637357401543597 },
    @{ Row = 5;  A = "04.jpg"; F = $null },
    @{ Row = 6;  A = "05.jpg"; F = 0.03033984992984335 },
    @{ Row = 7;  A = "06.jpg"; F = 0.0223175500808994 },
    @{ Row = 8;  A = "07.jpg"; F = 0.03061937810892377 },
    @{ Row = 9;  A = "08.jpg"; F = $null },
    @{ Row = 10; A = "09.jpg"; F = $null },
    @{ Row = 11; A = "10.jpg"; F = 0.03024929369681965 },
    @{ Row = 12; A = "11.jpg"; F = $null },
    @{ Row = 13; A = "12.jpg"; F = 0.04358519734800191 },
    @{ Row = 14; A = "13.jpg"; F = 0.02889985840097764 },
    @{ Row = 15; A = "14.jpg"; F = $null },
    @{ Row = 16; A = "15.jpg"; F = 0.02277738808077557 },
    @{ Row = 17; A = "16.jpg"; F = 0.021025794516685 },
    @{ Row = 18; A = "17.jpg"; F = 0.02189514940431469 },
    @{ Row = 19; A = "18.jpg"; F = $null },
    @{ Row = 20; A = "19.jpg"; F = $null },
    @{ Row = 21; A = "20.jpg"; F = 0.02408844832033849 },
    @{ Row = 22; A = "21.jpg"; F = $null },
    @{ Row = 23; A = "22.jpg"; F = $null },
    @{ Row = 24; A = "23.jpg"; F = $null },
    @{ Row = 25; A = "24.jpg"; F = 0.01952270447514764 },
    @{ Row = 26; A = "25.jpg"; F = $null },
    @{ Row = 27; A = "26.jpg"; F = 0.04119880719597929 },
    @{ Row = 28; A = "27.jpg"; F = 0.02340937419056989 },
    @{ Row = 29; A = "28.jpg"; F = $null },
    @{ Row = 30; A = "29.jpg"; F = $null },
    @{ Row = 31; A = "30.jpg"; F = 0.0304184855363709 },
    @{ Row = 32; A = "31.jpg"; F = $null },
    @{ Row = 33; A = "32.jpg"; F = 0.03572362339409182 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A
    if ($null -ne $u.F) {
        $ws.Cells.Item($r, 6).Value = $u.F
    }
}
